# Weekly update: insert this week's 3 new "Chirimoya" price rows
# (Terminal La Palmera de La Serena, Coquimbo, Provincia de Limari)
# above the historical data, shifting the existing rows down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at row 248 (pushes old 248.. down to 251..)
$ws.Rows.Item(248).Insert()
$ws.Rows.Item(248).Insert()
$ws.Rows.Item(248).Insert()

# Columns A,B,C,E,F,G,H,I,J,K are constant for every "Chirimoya" row in this
# block, so copy them straight from the (now shifted) row 251 for the new
# rows 248-250.
$srcRow = 251
$newRows = @(248, 249, 250)

foreach ($r in $newRows) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($srcRow, 1).Value()   # A Mercado ID
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($srcRow, 2).Value()   # B Mercado
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($srcRow, 3).Value()   # C Region
    $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($srcRow, 5).Value()   # E Codreg
    $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($srcRow, 6).Value()   # F Tipo
    $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($srcRow, 7).Value()   # G Producto ID
    $ws.Cells.Item($r, 8).Value = $ws.Cells.Item($srcRow, 8).Value()   # H Producto
    $ws.Cells.Item($r, 9).Value = $ws.Cells.Item($srcRow, 9).Value()   # I Categoria ID
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($srcRow, 10).Value() # J Categoria
    $ws.Cells.Item($r, 11).Value = $ws.Cells.Item($srcRow, 11).Value() # K Variedad
}

# Row 248: Chirimoya "Especial", $/bandeja 10 kilos
$ws.Cells.Item(248, 4).Value  = 45166
$ws.Cells.Item(248, 12).Value = "Especial"
$ws.Cells.Item(248, 13).Value = 160
$ws.Cells.Item(248, 14).Value = 24000
$ws.Cells.Item(248, 15).Value = 25000
$ws.Cells.Item(248, 16).Value = 24500
$ws.Cells.Item(248, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(248, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(248, 19).Value = 2450
$ws.Cells.Item(248, 20).Value = 10

# Row 249: Chirimoya "Primera", $/bandeja 10 kilos
$ws.Cells.Item(249, 4).Value  = 45166
$ws.Cells.Item(249, 12).Value = "Primera"
$ws.Cells.Item(249, 13).Value = 240
$ws.Cells.Item(249, 14).Value = 22000
$ws.Cells.Item(249, 15).Value = 23000
$ws.Cells.Item(249, 16).Value = 22500
$ws.Cells.Item(249, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(249, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(249, 19).Value = 2250
$ws.Cells.Item(249, 20).Value = 10

# Row 250: Chirimoya "Segunda", $/bandeja 10 kilos
$ws.Cells.Item(250, 4).Value  = 45166
$ws.Cells.Item(250, 12).Value = "Segunda"
$ws.Cells.Item(250, 13).Value = 240
$ws.Cells.Item(250, 14).Value = 18000
$ws.Cells.Item(250, 15).Value = 19000
$ws.Cells.Item(250, 16).Value = 18500
$ws.Cells.Item(250, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(250, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(250, 19).Value = 1850
$ws.Cells.Item(250, 20).Value = 10
